# The commit adds one new weekly price-record row for "Uva" (Terminal
# Hortofrutícola Agro Chillán) right before the existing row 52, which
# pushes all following records (old rows 52-155) down by one row
# (new rows 53-156). This mirrors a standard "insert a new observation at
# the top of the dated log" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; existing rows 52-155 shift down to 53-156.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record's data.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44953
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100109
$ws.Range("H52").Value = "Uva"
$ws.Range("I52").Value = 100109001
$ws.Range("J52").Value = "Uva"
$ws.Range("K52").Value = "Superior Seedless"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 60
$ws.Range("N52").Value = 8000
$ws.Range("O52").Value = 8000
$ws.Range("P52").Value = 8000
$ws.Range("Q52").Value = "`$/bandeja 10 kilos"
$ws.Range("R52").Value = "Provincia de Limarí"
$ws.Range("S52").Value = 800
$ws.Range("T52").Value = 10

Write-Host ("Row52 D (serial): " + $ws.Range("D52").Value2())
Write-Host ("Row52 K: " + $ws.Range("K52").Value())
Write-Host ("Row53 D (serial, shifted from old row52): " + $ws.Range("D53").Value2())
Write-Host ("Row156 K (shifted from old row155): " + $ws.Range("K156").Value())
Write-Host ("UsedRange rows: " + $ws.UsedRange.Rows.Count)
